$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two header/title cells
$ws.Range("B1").Value = "Example TableOld"
$ws.Range("B2").Value = "Example TableNew"

# Copy the style of row 1 (which uniformly uses the light-red highlight
# style, s="1", across all columns A-G) onto the two new trailing rows.
$ws.Range("A1:G1").Copy() | Out-Null
$ws.Range("A13:G13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A14:G14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 13: just a "-" marker in column A, rest blank
$ws.Range("A13").Value = "-"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = ""

# Row 14: "-" marker in column A plus trailing-row text in column B
$ws.Range("A14").Value = "-"
$ws.Range("B14").Value = "Trailing row here"
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("G14").Value = ""
